$wb = $excel.ActiveWorkbook

$oldId = "431d1843-bd14-45d9-9387-4015ca2b9a76"
$newId = "29681f68-d159-430e-91ca-adff909ec41c"

$oldHash = "400dcc228595f326ad3b27ed963e322bd1ab34a4"
$newHash = "5bff856cbfc9744e4211c06b72ff0cbbeab26935"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Hyperlinks.Delete()

$ws1.Range("A2").Value2 = "$newId.md"
$ws1.Range("B2").Value2 = "e2e\$newId.md"
$ws1.Range("G2").Value2 = "2016-08-27 02:58:31"

$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e0a60aeea2df2b8938581f484acb88cbcfc2e8a/e2e/$newId.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Drop every hyperlink on the sheet (engine quirk: Range.Hyperlinks.Delete()
# clears the whole sheet regardless of the range it was called on), then
# rebuild only the one we want to keep (A2). I2's link is intentionally not
# re-created.
$ws2.Range("A2").Hyperlinks.Delete()

$ws2.Range("G2").Value2 = "$newId.$newHash.zh-cn.xlf"
$ws2.Range("H2").Value2 = "2016-08-27 02:58:24"
$ws2.Range("I2").Value2 = ""
$ws2.Range("J2").Value2 = ""
$ws2.Range("K2").Value2 = "0001-01-01 00:00:00"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e0a60aeea2df2b8938581f484acb88cbcfc2e8a/e2e/$newId.md") | Out-Null

$ws2.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws2.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Hyperlinks.Delete()

$ws3.Range("G2").Value2 = "$newId.$newHash.de-de.xlf"
$ws3.Range("I2").Value2 = ""
$ws3.Range("J2").Value2 = ""
$ws3.Range("K2").Value2 = "0001-01-01 00:00:00"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8e0a60aeea2df2b8938581f484acb88cbcfc2e8a/e2e/$newId.md") | Out-Null

$ws3.Columns.Item(9).ColumnWidth = 18.6506053379604
$ws3.Columns.Item(10).ColumnWidth = 21.7054770333426
